$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new training/attendance date is appended right after the current last
# date column (BY, 2025-11-12 / serial 45973) as column BZ
# (2025-11-13 / serial 45974). For every row that already has data in BY,
# the same attendance value is copied into the new BZ column. Row 12's
# attendance data stops earlier in the sheet (no BY12 cell at all), so it is
# intentionally skipped and gets no BZ12 cell either -- matching the source
# edit exactly. Row 21 has no attendance value yet (BY21 is an empty,
# styled cell), so BZ21 is added the same way: styled but empty.
$rowValues = @{
    1  = 45974
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "B"
    6  = "B"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "B"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = $null
    22 = "P"
    23 = "P"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "RH"
}

foreach ($r in $rowValues.Keys) {
    $val = $rowValues[$r]
    if ($null -ne $val) {
        # Write the value first so dependent COUNTA/COUNTIF totals recalc.
        $ws.Range("BZ$r").Value = $val
    }
    # Then bring over BY's formatting (date/style) so BZ matches its neighbour.
    $ws.Range("BY$r").Copy() | Out-Null
    $ws.Range("BZ$r").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# Move the active selection to reflect where the user was last working.
$ws.Range("CB21").Select() | Out-Null
